$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 'b'
$ws.Range("J2").Value = 'Acknowledge (Backchannel)'
$ws.Range("I6").Value = 'sv'
$ws.Range("J6").Value = 'Statement-opinion'
$ws.Range("I8").Value = 'sd'
$ws.Range("J8").Value = 'Statement-non-opinion'
$ws.Range("I22").Value = 'b'
$ws.Range("J22").Value = 'Acknowledge (Backchannel)'
$ws.Range("I27").Value = 'sd'
$ws.Range("J27").Value = 'Statement-non-opinion'
$ws.Range("I31").Value = 'sd'
$ws.Range("J31").Value = 'Statement-non-opinion'
$ws.Range("I35").Value = 'sv'
$ws.Range("J35").Value = 'Statement-opinion'
$ws.Range("I44").Value = 'sd'
$ws.Range("J44").Value = 'Statement-non-opinion'
$ws.Range("I50").Value = 'b'
$ws.Range("J50").Value = 'Acknowledge (Backchannel)'
$ws.Range("I66").Value = 'sd'
$ws.Range("J66").Value = 'Statement-non-opinion'
$ws.Range("I71").Value = 'aa'
$ws.Range("J71").Value = 'Agree/Accept'
$ws.Range("I72").Value = 'sd'
$ws.Range("J72").Value = 'Statement-non-opinion'
$ws.Range("I73").Value = 'sd'
$ws.Range("J73").Value = 'Statement-non-opinion'
$ws.Range("I84").Value = 'aa'
$ws.Range("J84").Value = 'Agree/Accept'
$ws.Range("I89").Value = 'sd'
$ws.Range("J89").Value = 'Statement-non-opinion'
$ws.Range("I100").Value = 'sd'
$ws.Range("J100").Value = 'Statement-non-opinion'
$ws.Range("I111").Value = 'b'
$ws.Range("J111").Value = 'Acknowledge (Backchannel)'
$ws.Range("I112").Value = 'sd'
$ws.Range("J112").Value = 'Statement-non-opinion'
$ws.Range("I120").Value = 'sd'
$ws.Range("J120").Value = 'Statement-non-opinion'
$ws.Range("I123").Value = '%'
$ws.Range("J123").Value = 'Uninterpretable'
$ws.Range("I130").Value = 'aa'
$ws.Range("J130").Value = 'Agree/Accept'
$ws.Range("I153").Value = 'sv'
$ws.Range("J153").Value = 'Statement-opinion'
$ws.Range("I154").Value = 'ba'
$ws.Range("J154").Value = 'Appreciation'
$ws.Range("I166").Value = 'qy'
$ws.Range("J166").Value = 'Yes-No-Question'
$ws.Range("I176").Value = 'b'
$ws.Range("J176").Value = 'Acknowledge (Backchannel)'
$ws.Range("I177").Value = 'sd'
$ws.Range("J177").Value = 'Statement-non-opinion'
$ws.Range("I191").Value = 'sd'
$ws.Range("J191").Value = 'Statement-non-opinion'
$ws.Range("I197").Value = 'sd'
$ws.Range("J197").Value = 'Statement-non-opinion'
$ws.Range("I198").Value = '%'
$ws.Range("J198").Value = 'Uninterpretable'
$ws.Range("I206").Value = 'sd'
$ws.Range("J206").Value = 'Statement-non-opinion'
$ws.Range("I208").Value = 'sv'
$ws.Range("J208").Value = 'Statement-opinion'
$ws.Range("I212").Value = 'ba'
$ws.Range("J212").Value = 'Appreciation'
$ws.Range("I215").Value = 'aa'
$ws.Range("J215").Value = 'Agree/Accept'
$ws.Range("I230").Value = 'sd'
$ws.Range("J230").Value = 'Statement-non-opinion'
$ws.Range("I235").Value = 'b'
$ws.Range("J235").Value = 'Acknowledge (Backchannel)'
$ws.Range("I243").Value = 'sv'
$ws.Range("J243").Value = 'Statement-opinion'
$ws.Range("I257").Value = 'sd'
$ws.Range("J257").Value = 'Statement-non-opinion'
$ws.Range("I263").Value = 'aa'
$ws.Range("J263").Value = 'Agree/Accept'
$ws.Range("I282").Value = 'sd'
$ws.Range("J282").Value = 'Statement-non-opinion'
$ws.Range("I283").Value = 'sd'
$ws.Range("J283").Value = 'Statement-non-opinion'
$ws.Range("I294").Value = 'sd'
$ws.Range("J294").Value = 'Statement-non-opinion'
$ws.Range("I317").Value = 'sd'
$ws.Range("J317").Value = 'Statement-non-opinion'
$ws.Range("I318").Value = 'sd'
$ws.Range("J318").Value = 'Statement-non-opinion'
$ws.Range("I330").Value = 'sv'
$ws.Range("J330").Value = 'Statement-opinion'
$ws.Range("I331").Value = 'sd'
$ws.Range("J331").Value = 'Statement-non-opinion'
$ws.Range("I363").Value = 'sd'
$ws.Range("J363").Value = 'Statement-non-opinion'
$ws.Range("I373").Value = 'aa'
$ws.Range("J373").Value = 'Agree/Accept'
$ws.Range("I380").Value = 'sd'
$ws.Range("J380").Value = 'Statement-non-opinion'
$ws.Range("I386").Value = 'sd'
$ws.Range("J386").Value = 'Statement-non-opinion'
$ws.Range("I398").Value = 'sd'
$ws.Range("J398").Value = 'Statement-non-opinion'
$ws.Range("I401").Value = '%'
$ws.Range("J401").Value = 'Uninterpretable'
$ws.Range("I407").Value = 'sd'
$ws.Range("J407").Value = 'Statement-non-opinion'
$ws.Range("I416").Value = 'sv'
$ws.Range("J416").Value = 'Statement-opinion'
$ws.Range("I443").Value = 'b'
$ws.Range("J443").Value = 'Acknowledge (Backchannel)'
$ws.Range("I444").Value = 'ba'
$ws.Range("J444").Value = 'Appreciation'
$ws.Range("I449").Value = 'sd'
$ws.Range("J449").Value = 'Statement-non-opinion'
$ws.Range("I477").Value = 'sv'
$ws.Range("J477").Value = 'Statement-opinion'
$ws.Range("I483").Value = 'sd'
$ws.Range("J483").Value = 'Statement-non-opinion'
$ws.Range("I510").Value = 'b'
$ws.Range("J510").Value = 'Acknowledge (Backchannel)'
$ws.Range("I515").Value = 'b'
$ws.Range("J515").Value = 'Acknowledge (Backchannel)'
$ws.Range("I527").Value = '%'
$ws.Range("J527").Value = 'Uninterpretable'
$ws.Range("I528").Value = '%'
$ws.Range("J528").Value = 'Uninterpretable'
$ws.Range("I537").Value = 'b'
$ws.Range("J537").Value = 'Acknowledge (Backchannel)'
$ws.Range("I539").Value = 'sd'
$ws.Range("J539").Value = 'Statement-non-opinion'
$ws.Range("I542").Value = 'ba'
$ws.Range("J542").Value = 'Appreciation'
$ws.Range("I550").Value = 'sd'
$ws.Range("J550").Value = 'Statement-non-opinion'
$ws.Range("I554").Value = 'aa'
$ws.Range("J554").Value = 'Agree/Accept'
$ws.Range("I556").Value = 'b'
$ws.Range("J556").Value = 'Acknowledge (Backchannel)'
$ws.Range("I562").Value = 'sd'
$ws.Range("J562").Value = 'Statement-non-opinion'
$ws.Range("I566").Value = 'sd'
$ws.Range("J566").Value = 'Statement-non-opinion'
$ws.Range("I573").Value = 'aa'
$ws.Range("J573").Value = 'Agree/Accept'
$ws.Range("I575").Value = 'sd'
$ws.Range("J575").Value = 'Statement-non-opinion'
$ws.Range("I582").Value = 'b'
$ws.Range("J582").Value = 'Acknowledge (Backchannel)'
$ws.Range("I596").Value = 'sv'
$ws.Range("J596").Value = 'Statement-opinion'
$ws.Range("I602").Value = 'b'
$ws.Range("J602").Value = 'Acknowledge (Backchannel)'
